$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0.001
$ws.Range("K12").Value = 475
$ws.Range("L12").Value = 0.002375
